$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-indexed, matching the original sheet) that need to be removed.
# Sorted descending so earlier deletions don't shift the indices of rows
# still pending deletion.
$rowsToDelete = @(184, 158, 149, 147, 133, 132, 131, 64, 62, 59, 44)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Row deletion in this runtime does not resynchronize the worksheet
# hyperlinks (they keep pointing at their original, now-stale rows), so
# rebuild the hyperlink collection for column C from scratch.
$ws.Cells.Hyperlinks.Delete()

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $url = $cell.Value2
    if ($url) {
        $ws.Hyperlinks.Add($cell, $url)
        $cell.Style = "Hyperlink"
    }
}
